$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").Value = "apr"
$ws.Range("D10").Value = "40K"

# Move the active selection to the next empty cell below the new row,
# matching Excel's natural behavior after data entry.
$ws.Range("D11").Select() | Out-Null
